$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 2686
$ws.Range("K3").Value = 2595
$ws.Range("E4").Value = 2032
$ws.Range("G4").Value = 1484
$ws.Range("H4").Value = 1727
$ws.Range("I4").Value = 1788
$ws.Range("K4").Value = 543
$ws.Range("K5").Value = 173
$ws.Range("K6").Value = 3234
$ws.Range("E7").Value = 26037
$ws.Range("G7").Value = 24710
$ws.Range("H7").Value = 26040
$ws.Range("I7").Value = 26241
$ws.Range("K7").Value = 9231

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K3").Value = 30
$ws.Range("K6").Value = 71
$ws.Range("K7").Value = 136

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K6").Value = 205
$ws.Range("K7").Value = 608

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K3").Value = 132
$ws.Range("K7").Value = 364

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 46
$ws.Range("K3").Value = 52
$ws.Range("K7").Value = 149

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K6").Value = 94
$ws.Range("K7").Value = 302

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 65
$ws.Range("K6").Value = 89
$ws.Range("K7").Value = 217

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 68
$ws.Range("K6").Value = 75
$ws.Range("K7").Value = 277
$ws.Range("K8").Value = 608
$ws.Range("K9").Value = 34
$ws.Range("K11").Value = 194
$ws.Range("K15").Value = 93
$ws.Range("K18").Value = 62
$ws.Range("K19").Value = 269
$ws.Range("K20").Value = 211
$ws.Range("K22").Value = 28
$ws.Range("K23").Value = 81
$ws.Range("K24").Value = 32
$ws.Range("K27").Value = 100
$ws.Range("K29").Value = 481
$ws.Range("K31").Value = 107
$ws.Range("K33").Value = 364
$ws.Range("K37").Value = 302
$ws.Range("K39").Value = 13
$ws.Range("K41").Value = 82
$ws.Range("K42").Value = 322
$ws.Range("K51").Value = 101
$ws.Range("K52").Value = 254
$ws.Range("K53").Value = 136
$ws.Range("K54").Value = 170
$ws.Range("K60").Value = 60
$ws.Range("E63").Value = 368
$ws.Range("G63").Value = 285
$ws.Range("H63").Value = 280
$ws.Range("I63").Value = 201
$ws.Range("K63").Value = 33
$ws.Range("K65").Value = 217
$ws.Range("K67").Value = 358
$ws.Range("K68").Value = 24
$ws.Range("K71").Value = 29
$ws.Range("K72").Value = 42
$ws.Range("K76").Value = 137
$ws.Range("K85").Value = 443
$ws.Range("K86").Value = 57
$ws.Range("K89").Value = 121
$ws.Range("K90").Value = 85
$ws.Range("K91").Value = 87
$ws.Range("K94").Value = 110
$ws.Range("K95").Value = 149
$ws.Range("E101").Value = 26037
$ws.Range("G101").Value = 24710
$ws.Range("H101").Value = 26040
$ws.Range("I101").Value = 26241
$ws.Range("K101").Value = 9231

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K6").Value = 43
$ws.Range("K7").Value = 107

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 115
$ws.Range("K3").Value = 112
$ws.Range("K7").Value = 358

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 76
$ws.Range("K7").Value = 170

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 128
$ws.Range("K3").Value = 160
$ws.Range("K6").Value = 155
$ws.Range("K7").Value = 481

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K4").Value = 10
$ws.Range("K7").Value = 269

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K3").Value = 24
$ws.Range("K7").Value = 137

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K3").Value = 25
$ws.Range("K7").Value = 75

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K2").Value = 29
$ws.Range("K7").Value = 82

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 82
$ws.Range("K3").Value = 102
$ws.Range("K4").Value = 12
$ws.Range("K6").Value = 124
$ws.Range("K7").Value = 322

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K2").Value = 28
$ws.Range("K7").Value = 81

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 74
$ws.Range("K6").Value = 73
$ws.Range("K7").Value = 211

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K4").Value = 8
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K6").Value = 76
$ws.Range("K7").Value = 277

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 30
$ws.Range("K3").Value = 21
$ws.Range("K6").Value = 48
$ws.Range("K7").Value = 110

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K3").Value = 22
$ws.Range("K4").Value = 5
$ws.Range("K6").Value = 31
$ws.Range("K7").Value = 93

$ws = $wb.Worksheets.Item("Greektown")
$ws.Range("K5").Value = 7
$ws.Range("K6").Value = 13

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K6").Value = 76
$ws.Range("K7").Value = 194

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K2").Value = 11
$ws.Range("K7").Value = 34

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K6").Value = 26
$ws.Range("K7").Value = 68

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K4").Value = 18
$ws.Range("K7").Value = 121

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K3").Value = 21
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 100

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 85

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 101

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("K3").Value = 7
$ws.Range("K7").Value = 24

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 161
$ws.Range("K3").Value = 152
$ws.Range("K5").Value = 10
$ws.Range("K7").Value = 443

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("K2").Value = 15
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("K2").Value = 13
$ws.Range("K7").Value = 29

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 42

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K3").Value = 60
$ws.Range("K6").Value = 104
$ws.Range("K7").Value = 254
